$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.021.20'
$ws.Range('E2').Value = '  -3.76%  '
$ws.Range('D3').Value = '3.225.07'
$ws.Range('E3').Value = '  -4.58%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '539.98'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.80%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.63'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -10.13%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.228.01'
$ws.Range('E8').Value = '  -4.40%  '
$ws.Range('E9').Value = '  -5.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.61'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.116'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -5.80%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.398'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -5.03%  '
$ws.Range('D13').Value = '3.778.44'
$ws.Range('E13').Value = '  -4.73%  '
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.15'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -7.93%  '
$ws.Range('D16').Value = '3.218.42'
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000160'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -6.72%  '
$ws.Range('D18').Value = '59.130.34'
$ws.Range('E18').Value = '  -3.84%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.93'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.86%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.25'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -7.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.31'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '362.20'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.53%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.46'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -6.73%  '
$ws.Range('E25').Value = '  -8.21%  '
$ws.Range('D26').Value = '3.355.07'
$ws.Range('E26').Value = '  -5.31%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0₃0979'
$ws.Range('E27').Value = '  -10.46%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.171'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.01'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.05'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.72%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  -8.11%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.09'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -8.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '21.97'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.64%  '
$ws.Range('E35').Value = '  -3.01%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.94'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -8.83%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '162.35'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.40'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.20%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.45'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -7.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '25.65'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -15.71%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0706'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -7.92%  '
$ws.Range('D42').Value = '3.257.01'
$ws.Range('E42').Value = '  -4.86%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '41.05'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('E44').Value = '  -6.01%  '
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('E46').Value = '  -7.90%  '
$ws.Range('E47').Value = '  -6.84%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = '2.289.29'
$ws.Range('E49').Value = '  -9.35%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.27'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -7.23%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '20.76'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -9.62%  '
